$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7170026666666667
$ws.Range("H2").Value = 2.151008
$ws.Range("I2").Value = 0.02953485643833859
$ws.Range("J2").Value = 0.02953485643833859
$ws.Range("M2").Value = 36.89194233333333
$ws.Range("N2").Value = 110.675827
$ws.Range("O2").Value = 0.3567095043190808
$ws.Range("P2").Value = 0.3567095043190809
$ws.Range("Q2").Value = 26.45162103151289
$ws.Range("R2").Value = 238.064589283616
$ws.Range("S2").Value = 0.01053536400025497
$ws.Range("T2").Value = 0.01053536400025497

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7170026666666667
$ws.Range("H3").Value = 2.151008
$ws.Range("I3").Value = 0.02953485643833859
$ws.Range("J3").Value = 0.02953485643833859
$ws.Range("M3").Value = 42.68037399999999
$ws.Range("N3").Value = 128.041122
$ws.Range("O3").Value = 0.4126780562577495
$ws.Range("P3").Value = 0.4126780562577496
$ws.Range("Q3").Value = 30.60194197233066
$ws.Range("R3").Value = 275.4174777509759
$ws.Range("S3").Value = 0.01218838714682525
$ws.Range("T3").Value = 0.01218838714682525

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7170026666666667
$ws.Range("H4").Value = 2.151008
$ws.Range("I4").Value = 0.02953485643833859
$ws.Range("J4").Value = 0.02953485643833859
$ws.Range("M4").Value = 23.85061433333334
$ws.Range("N4").Value = 71.55184300000001
$ws.Range("O4").Value = 0.2306124394231696
$ws.Range("P4").Value = 0.2306124394231696
$ws.Range("Q4").Value = 17.10095407863822
$ws.Range("R4").Value = 153.908586707744
$ws.Range("S4").Value = 0.006811105291258368
$ws.Range("T4").Value = 0.006811105291258369

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.4970672037825566
$ws.Range("J5").Value = 0.4970672037825566
$ws.Range("M5").Value = 36.89194233333333
$ws.Range("N5").Value = 110.675827
$ws.Range("O5").Value = 0.3567095043190808
$ws.Range("P5").Value = 0.3567095043190809
$ws.Range("Q5").Value = 445.1768143549371
$ws.Range("R5").Value = 4006.591329194433
$ws.Range("S5").Value = 0.1773085958745473
$ws.Range("T5").Value = 0.1773085958745473

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.4970672037825566
$ws.Range("J6").Value = 0.4970672037825566
$ws.Range("M6").Value = 42.68037399999999
$ws.Range("N6").Value = 128.041122
$ws.Range("O6").Value = 0.4126780562577495
$ws.Range("P6").Value = 0.4126780562577496
$ws.Range("Q6").Value = 515.0260932623693
$ws.Range("R6").Value = 4635.234839361323
$ws.Range("S6").Value = 0.2051287274864601
$ws.Range("T6").Value = 0.2051287274864602

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.4970672037825566
$ws.Range("J7").Value = 0.4970672037825566
$ws.Range("M7").Value = 23.85061433333334
$ws.Range("N7").Value = 71.55184300000001
$ws.Range("O7").Value = 0.2306124394231696
$ws.Range("P7").Value = 0.2306124394231696
$ws.Range("Q7").Value = 287.8064920894118
$ws.Range("S7").Value = 0.1146298804215491
$ws.Range("T7").Value = 0.1146298804215491

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.4733979397791048
$ws.Range("J8").Value = 0.4733979397791048
$ws.Range("M8").Value = 36.89194233333333
$ws.Range("N8").Value = 110.675827
$ws.Range("O8").Value = 0.3567095043190808
$ws.Range("P8").Value = 0.3567095043190809
$ws.Range("Q8").Value = 423.9784583439215
$ws.Range("R8").Value = 3815.806125095294
$ws.Range("S8").Value = 0.1688655444442786
$ws.Range("T8").Value = 0.1688655444442786

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.4733979397791048
$ws.Range("J9").Value = 0.4733979397791048
$ws.Range("M9").Value = 42.68037399999999
$ws.Range("N9").Value = 128.041122
$ws.Range("O9").Value = 0.4126780562577495
$ws.Range("P9").Value = 0.4126780562577496
$ws.Range("Q9").Value = 490.5016658261426
$ws.Range("R9").Value = 4414.514992435284
$ws.Range("S9").Value = 0.1953609416244641
$ws.Range("T9").Value = 0.1953609416244642

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.4733979397791048
$ws.Range("J10").Value = 0.4733979397791048
$ws.Range("M10").Value = 23.85061433333334
$ws.Range("N10").Value = 71.55184300000001
$ws.Range("O10").Value = 0.2306124394231696
$ws.Range("P10").Value = 0.2306124394231696
$ws.Range("S10").Value = 0.1091714537103621
$ws.Range("T10").Value = 0.1091714537103621
